$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update age/group header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: clear deleted values, update remaining one
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -0.86568835424872725

# Row 3: clear deleted value, update/add remaining ones
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 0.11022866506253015
$ws.Range("D3").Value = -1.7914419616663402
$ws.Range("E3").Value = 2.0033802731607646

# Update selection to match new narrower highlighted range
$ws.Range("B1:E3").Select()
